$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 53
$ws.Range("H53").Value = 142.85715
$ws.Range("I53").Value = 190.2
$ws.Range("K53").Value = 190.2
$ws.Range("M53").Value = 446.8
# Row 64
$ws.Range("H64").Value = 5779.8
$ws.Range("I64").Value = 6424.75
$ws.Range("K64").Value = 6424.75
$ws.Range("M64").Value = -6176.75
# Row 67
$ws.Range("H67").Value = 5779.8
$ws.Range("I67").Value = 6424.75
$ws.Range("K67").Value = 6424.75
$ws.Range("M67").Value = -5566.75
# Row 88
$ws.Range("H88").Value = 6000
$ws.Range("I88").Value = 5000
$ws.Range("J88").Value = 6500
$ws.Range("K88").Value = 5000
$ws.Range("L88").Value = 6500
$ws.Range("M88").Value = -4594
$ws.Range("N88").Value = -7312
# Row 91
$ws.Range("H91").Value = 6000
$ws.Range("I91").Value = 5000
$ws.Range("J91").Value = 6500
$ws.Range("K91").Value = 5000
$ws.Range("L91").Value = 6500
$ws.Range("M91").Value = -3596
$ws.Range("N91").Value = -9308
# Row 137
$ws.Range("H137").Value = 1696.3334
$ws.Range("I137").Value = 1474.5
$ws.Range("J137").Value = 2140
$ws.Range("K137").Value = 4423.5
$ws.Range("L137").Value = 6420
$ws.Range("M137").Value = -1873.5
$ws.Range("N137").Value = -11520

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 50.833332
$ws.Range("I5").Value = 43
$ws.Range("J5").Value = 90
$ws.Range("K5").Value = 43
$ws.Range("L5").Value = 90
$ws.Range("M5").Value = 69
$ws.Range("N5").Value = -314
# Row 32
$ws.Range("H32").Value = 4101.967
$ws.Range("I32").Value = 1918.7084
$ws.Range("K32").Value = 1918.7084
$ws.Range("M32").Value = -1631.7084
# Row 45
$ws.Range("H45").Value = 2501.375
$ws.Range("J45").Value = 3300
$ws.Range("L45").Value = 3300
$ws.Range("N45").Value = -4054
# Row 61
$ws.Range("H61").Value = 4616.7856
$ws.Range("I61").Value = 4725.769
$ws.Range("J61").Value = 3200
$ws.Range("K61").Value = 4725.769
$ws.Range("L61").Value = 3200
$ws.Range("M61").Value = -4513.769
$ws.Range("N61").Value = -3624
# Row 122
$ws.Range("H122").Value = 444709
$ws.Range("I122").Value = 532020.44
$ws.Range("K122").Value = 1596061.32
$ws.Range("M122").Value = -1593611.32
# Row 132
$ws.Range("H132").Value = 1723.5454
$ws.Range("I132").Value = 1695
$ws.Range("J132").Value = 2009
$ws.Range("K132").Value = 5085
$ws.Range("L132").Value = 6027
$ws.Range("M132").Value = -2555
$ws.Range("N132").Value = -11087
# Row 136
$ws.Range("H136").Value = 4616.7856
$ws.Range("I136").Value = 4725.769
$ws.Range("J136").Value = 3200
$ws.Range("K136").Value = 14177.307
$ws.Range("L136").Value = 9600
$ws.Range("M136").Value = -11627.307
$ws.Range("N136").Value = -14700

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 50.833332
$ws.Range("I4").Value = 43
$ws.Range("J4").Value = 90
$ws.Range("K4").Value = 43
$ws.Range("L4").Value = 90
$ws.Range("M4").Value = 72
$ws.Range("N4").Value = -320
# Row 14
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
# Row 60
$ws.Range("H60").Value = 49332.332
$ws.Range("J60").Value = 49332.332
$ws.Range("L60").Value = 49332.332
$ws.Range("N60").Value = -50530.332

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 4524.75
$ws.Range("I16").Value = 4439.6
$ws.Range("J16").Value = 4666.6665
$ws.Range("K16").Value = 4439.6
$ws.Range("L16").Value = 4666.6665
$ws.Range("M16").Value = -4152.6
$ws.Range("N16").Value = -5240.6665
# Row 62
$ws.Range("H62").Value = 197254.5
$ws.Range("I62").Value = 129673
$ws.Range("J62").Value = 399999
$ws.Range("K62").Value = 129673
$ws.Range("L62").Value = 399999
$ws.Range("M62").Value = -129049
$ws.Range("N62").Value = -401247
# Row 65
$ws.Range("H65").Value = 197254.5
$ws.Range("I65").Value = 129673
$ws.Range("J65").Value = 399999
$ws.Range("K65").Value = 648365
$ws.Range("L65").Value = 1999995
$ws.Range("M65").Value = -645245
$ws.Range("N65").Value = -2006235
# Row 105
$ws.Range("H105").Value = 2092.5334
$ws.Range("J105").Value = 3114.2856
$ws.Range("L105").Value = 3114.2856
$ws.Range("N105").Value = -6608.2856
# Row 113
$ws.Range("H113").Value = 4524.75
$ws.Range("I113").Value = 4439.6
$ws.Range("J113").Value = 4666.6665
$ws.Range("K113").Value = 4439.6
$ws.Range("L113").Value = 4666.6665
$ws.Range("M113").Value = -2269.6
$ws.Range("N113").Value = -9006.6665
# Row 125
$ws.Range("H125").Value = 89498
$ws.Range("J125").Value = 89498
$ws.Range("L125").Value = 89498
$ws.Range("N125").Value = -94418

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 24290814
$ws.Range("I4").Value = 32909228
$ws.Range("K4").Value = 98727684
$ws.Range("M4").Value = -98727572
# Row 23
$ws.Range("H23").Value = 250034.75
$ws.Range("J23").Value = 500049.5
$ws.Range("L23").Value = 1500148.5
$ws.Range("N23").Value = -1500618.5
# Row 57
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
# Row 98
$ws.Range("H98").Value = 3117.8
$ws.Range("J98").Value = 2772.375
$ws.Range("L98").Value = 8317.125
$ws.Range("N98").Value = -11313.125
# Row 113
$ws.Range("H113").Value = 723.75
$ws.Range("I113").Value = 550
$ws.Range("J113").Value = 897.5
$ws.Range("K113").Value = 1650
$ws.Range("L113").Value = 2692.5
$ws.Range("M113").Value = 520
$ws.Range("N113").Value = -7032.5
# Row 115
$ws.Range("H115").Value = 394
$ws.Range("I115").Value = 394
$ws.Range("K115").Value = 1182
$ws.Range("M115").Value = -7
# Row 131
$ws.Range("H131").Value = 1436.4572
$ws.Range("I131").Value = 816
$ws.Range("J131").Value = 1539.8667
$ws.Range("K131").Value = 2448
$ws.Range("L131").Value = 4619.6001
$ws.Range("M131").Value = 2592
$ws.Range("N131").Value = -14699.6001

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 35
$ws.Range("H35").Value = 3760000
$ws.Range("J35").Value = 5000000
$ws.Range("L35").Value = 5000000
$ws.Range("N35").Value = -5000596
# Row 80
$ws.Range("H80").Value = 4166.5557
$ws.Range("J80").Value = 4685.5713
$ws.Range("L80").Value = 4685.5713
$ws.Range("N80").Value = -6681.5713
# Row 83
$ws.Range("H83").Value = 4166.5557
$ws.Range("J83").Value = 4685.5713
$ws.Range("L83").Value = 23427.8565
$ws.Range("N83").Value = -33411.85649999999
# Row 97
$ws.Range("H97").Value = 2340.6
$ws.Range("I97").Value = 1822.3572
$ws.Range("J97").Value = 3000.182
$ws.Range("K97").Value = 1822.3572
$ws.Range("L97").Value = 3000.182
$ws.Range("M97").Value = -1326.3572
$ws.Range("N97").Value = -3992.182
# Row 102
$ws.Range("H102").Value = 2807.6
$ws.Range("I102").Value = 2374
$ws.Range("K102").Value = 2374
$ws.Range("M102").Value = -752
# Row 107
$ws.Range("H107").Value = 1222.8334
$ws.Range("I107").Value = 1293.6
$ws.Range("K107").Value = 1293.6
$ws.Range("M107").Value = 626.4000000000001
# Row 122
$ws.Range("H122").Value = 61105.53
$ws.Range("I122").Value = 1868.9231
$ws.Range("K122").Value = 5606.7693
$ws.Range("M122").Value = -3156.7693

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 1460.2307
$ws.Range("I61").Value = 867.1818
$ws.Range("J61").Value = 4722
$ws.Range("K61").Value = 867.1818
$ws.Range("L61").Value = 4722
$ws.Range("M61").Value = -665.1818
$ws.Range("N61").Value = -5126
# Row 82
$ws.Range("H82").Value = 15750
$ws.Range("I82").Value = 15750
$ws.Range("K82").Value = 15750
$ws.Range("M82").Value = -15389
# Row 85
$ws.Range("H85").Value = 15750
$ws.Range("I85").Value = 15750
$ws.Range("K85").Value = 15750
$ws.Range("M85").Value = -14502
# Row 113
$ws.Range("H113").Value = 1460.2307
$ws.Range("I113").Value = 867.1818
$ws.Range("J113").Value = 4722
$ws.Range("K113").Value = 867.1818
$ws.Range("L113").Value = 4722
$ws.Range("M113").Value = 1302.8182
$ws.Range("N113").Value = -9062

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 459
$ws.Range("I107").Value = 459
$ws.Range("K107").Value = 1377
$ws.Range("M107").Value = 543
# Row 124
$ws.Range("H124").Value = 79999.5
$ws.Range("J124").Value = 79999.5
$ws.Range("L124").Value = 79999.5
$ws.Range("N124").Value = -89819.5

Write-Host "Sheets updated successfully"
